$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "usuarioAp" column (column D), shifting E:G left to D:F.
$ws.Range("D1:D2").EntireColumn.Delete()

# Update the now-shifted C1 header text from "cuenta prestamo senior" to "cuenta prestamo".
$ws.Range("C1").Value = "cuenta prestamo"

# Restore the selection position recorded in the saved file.
$ws.Range("E5").Select()
